$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "dft2_z_s1_type1_no_op"
$ws.Range("D2").Value = "dft2_z_s1_type1"

$ws.Range("A3").Value = "dft2_z_s3_type1_no_op"
$ws.Range("D3").Value = "dft2_z_s3_type1"
